$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.968.41"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.23%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.765.45"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.26%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "646.61"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.84%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.49"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.29%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.764.35"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -1.22%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  +0.55%  "
$ws.Range("E10").Value = "  -2.23%  "
$ws.Range("E11").Value = "  +0.52%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.89"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +4.34%  "
$ws.Range("E13").Value = "  -4.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.83"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -3.12%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.400.25"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.16%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.758.25"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.57%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.929.58"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.29%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.66"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.64%  "
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("E20").Value = "  -1.68%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "471.66"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.55"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.10%  "
$ws.Range("E23").Value = "  -0.61%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000144"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -4.69%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.88"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.14"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.26%  "
$ws.Range("E27").Value = "  -2.71%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.07"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.914.55"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.21%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.68"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.54%  "
$ws.Range("E32").Value = "  +1.85%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.11"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.27%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "28.50"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.23%  "
$ws.Range("E35").Value = "  +17.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.720.44"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.82"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.65%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.77"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.15%  "
$ws.Range("E41").Value = "  -7.06%  "
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.956"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "44.87"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +4.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.97"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +3.50%  "
$ws.Range("E47").Value = "  -1.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "47.52"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.30%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.41"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.81%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.295"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.69%  "
$ws.Range("E51").Value = "  -0.97%  "
